# Weekly refresh of the "Perejil" (parsley) price table for Vega Modelo de
# Temuco: a new daily record is inserted right before the current row 45
# (shifting every subsequent record down by one row), and the new record's
# own values are written into the freshly inserted row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 45; Excel shifts rows 45:180 down to
# 46:181 and grows the used range / <dimension> to A1:R181 automatically.
$ws.Rows.Item(45).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A45").Value = 10
$ws.Range("B45").Value = "Vega Modelo de Temuco"
$ws.Range("C45").Value = "La Araucanía"
$ws.Range("D45").Value = 44453
$ws.Range("E45").Value = 9
$ws.Range("F45").Value = 100112044
$ws.Range("G45").Value = "Perejil"
$ws.Range("H45").Value = "Sin especificar"
$ws.Range("I45").Value = "Primera"
$ws.Range("J45").Value = 20
$ws.Range("K45").Value = 4000
$ws.Range("L45").Value = 4000
$ws.Range("M45").Value = 4000
$ws.Range("N45").Value = '$/docena de atados (3 kilos)'
$ws.Range("O45").Value = "Provincia de Cautín"
$ws.Range("P45").Value = 1333
$ws.Range("Q45").Value = 3
$ws.Range("R45").Value = "Hortaliza"
